$d = $word.ActiveDocument

# Paragraph 2 ("人生在世，匆匆一晃几十年就过去了，有些事是将就不了的") loses its
# list-paragraph style/numbering/indent and its run text, becoming an empty
# paragraph that only keeps the paragraph-mark run properties (rFonts hint,
# sz, szCs) that the old run used to carry.
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:rFonts w:hint=""eastAsia""/><w:sz w:val=""18""/><w:szCs w:val=""18""/></w:rPr></w:pPr></w:p>")

# Paragraphs 3 ("人张了嘴就是要把误会说清楚...浪费光阴! ") and 4 ("其实你也能
# 像风筝一样...不在落难时拂袖。") are removed entirely, marks and all.
$p3 = $d.Paragraphs.Item(3)
$p4 = $d.Paragraphs.Item(4)
$killRange = $d.Range($p3.Range.Start, $p4.Range.End)
$killRange.Delete()
